# Update "想去人数" (want-to-go count) figures scraped for the
# gh-pages data refresh (commit 456a3b4).
#
# Sheet 展览 (Exhibitions)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 275
$ws.Range("F3").Value = 947
$ws.Range("F7").Value = 145
$ws.Range("F8").Value = 1222
$ws.Range("F10").Value = 3143
$ws.Range("F14").Value = 659
$ws.Range("F15").Value = 29
$ws.Range("F16").Value = 541
$ws.Range("F17").Value = 266
$ws.Range("F20").Value = 1279
$ws.Range("F21").Value = 1279
$ws.Range("F22").Value = 206
$ws.Range("F25").Value = 239
$ws.Range("F27").Value = 404
$ws.Range("F30").Value = 279
$ws.Range("F32").Value = 870
$ws.Range("F33").Value = 147
$ws.Range("F35").Value = 355
$ws.Range("F37").Value = 5153
$ws.Range("F40").Value = 223
$ws.Range("F43").Value = 22

# Sheet 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 1783
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 53
$ws.Range("F17").Value = 314
$ws.Range("F23").Value = 751
$ws.Range("F26").Value = 17
$ws.Range("F42").Value = 16

# Sheet 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 464

# Sheet 全部类型 (All types, aggregated view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 275
$ws.Range("F5").Value = 947
$ws.Range("F9").Value = 145
$ws.Range("F10").Value = 1222
$ws.Range("F13").Value = 3143
$ws.Range("F18").Value = 659
$ws.Range("F19").Value = 464
$ws.Range("F20").Value = 29
$ws.Range("F21").Value = 541
$ws.Range("F22").Value = 266
$ws.Range("F26").Value = 1279
$ws.Range("F27").Value = 1279
$ws.Range("F28").Value = 206
$ws.Range("F29").Value = 53
$ws.Range("F31").Value = 239
$ws.Range("F32").Value = 404
$ws.Range("F36").Value = 279
$ws.Range("F37").Value = 870
$ws.Range("F38").Value = 147
$ws.Range("F40").Value = 355
$ws.Range("F42").Value = 5153
$ws.Range("F46").Value = 223
$ws.Range("F49").Value = 22
